$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gains two
#    trailing spaces, then a new parenthetical remark in red (C00000)
#    split across three runs.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$p1 = $d.Paragraphs.Item(1)
$insertPoint = $p1.Range.End - 1

$enDash = [char]0x2013
$part1 = "(This is a change " + $enDash + " Ve"
$part2 = "rsion for branch alternate"
$part3 = ")"

$r = $d.Range($insertPoint, $insertPoint)
$r.InsertAfter($part1)
$seg1 = $d.Range($insertPoint, $insertPoint + $part1.Length)
$seg1.Font.Color = 192

$pos2 = $insertPoint + $part1.Length
$r = $d.Range($pos2, $pos2)
$r.InsertAfter($part2)
$seg2 = $d.Range($pos2, $pos2 + $part2.Length)
$seg2.Font.Color = 192

$pos3 = $pos2 + $part2.Length
$r = $d.Range($pos3, $pos3)
$r.InsertAfter($part3)
$seg3 = $d.Range($pos3, $pos3 + $part3.Length)
$seg3.Font.Color = 192

# ---------------------------------------------------------------------
# 2) Append a new, otherwise-empty paragraph after the final paragraph
#    of the speech, shaded with fill F9F9F9. Building it via InsertXML
#    keeps it free of any inherited run/paragraph formatting.
# ---------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$endRange.InsertXML($newParaXml)

Write-Output "done"
